# Actualización automática de tasas-transfi.xlsx
# Updates the daily conversion note on "Hoja1" and the corresponding
# rate cells (N10/O10/N12) on the "tasas" sheet.

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion-of-the-day note (cell A1) ---
$hoja1 = $wb.Worksheets.Item("Hoja1")

$newNote = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 13.93 = 56026.74 pesos
✅ 56026.74 pesos = 13.9 = 963.38 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$hoja1.Range("A1").Value = $newNote.TrimEnd("`r", "`n")

# --- tasas: update the rate cells ---
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 71.80500000000001
$tasas.Range("O10").Value = 4023
$tasas.Range("N12").Value = 4032
